$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# --- Row 2 (users table field list): drop the "roles" field, shift remaining left ---
$ws.Range("E2").Value = "groups"
$ws.Range("F2").Value = "first_name"
$ws.Range("G2").Value = "last_name"
$ws.Range("H2").Value = ""

# --- Row 3 (users table field notes): drop the "23 ky tu" note, keep the "6 ky tu" note ---
$ws.Range("E3").Value = "6 ký tự`nXem thêm`ntrong  types.php"
$ws.Range("F3").Value = ""

# --- Row 18/19 (groups table): add the "roles" field + its note ---
$ws.Range("D18").Value = "roles"
$ws.Range("D19").Value = "23 ký tự`nXem thêm`ntrong  types.php"
$ws.Range("D19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 60

# --- Update sheet view: scroll + selection ---
$ws.Activate()
$ws.Range("E19").Select()
$excel.ActiveWindow.ScrollRow = 16
